$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Hunk 0 (diff @@ -1468,25 +1468,25 @@)
$ws.Range("H17").Value = 1212.1875
$ws.Range("J17").Value = 1254.3572
$ws.Range("L17").Value = 3763.0716
$ws.Range("N17").Value = -4099.071599999999

# Hunk 1 (diff @@ -4585,25 +4585,25 @@)
$ws.Range("H80").Value = 50594.6
$ws.Range("J80").Value = 371
$ws.Range("L80").Value = 1113
$ws.Range("N80").Value = -3109

# Hunk 2 (diff @@ -4732,25 +4732,25 @@)
$ws.Range("H83").Value = 50594.6
$ws.Range("J83").Value = 371
$ws.Range("L83").Value = 3339
$ws.Range("N83").Value = -13323

# Hunk 3 (diff @@ -5892,22 +5892,22 @@)
$ws.Range("H106").Value = 8093.5835
$ws.Range("I106").Value = 10143.286
$ws.Range("K106").Value = 10143.286
$ws.Range("M106").Value = -9512.286

# Hunk 4 (diff @@ -6847,22 +6847,22 @@)
$ws.Range("H125").Value = 7762.75
$ws.Range("I125").Value = 13035.667
$ws.Range("K125").Value = 117321.003
$ws.Range("M125").Value = -114861.003

# Hunk 5 (diff @@ -7438,25 +7438,25 @@)
$ws.Range("H137").Value = 7917.457
$ws.Range("I137").Value = 12148
$ws.Range("J137").Value = 2276.7334
$ws.Range("K137").Value = 36444
$ws.Range("L137").Value = 6830.2002
$ws.Range("M137").Value = -33894
$ws.Range("N137").Value = -11930.2002

# Hunk 6 (diff @@ -7640,25 +7640,25 @@)
$ws.Range("H141").Value = 6287.1113
$ws.Range("I141").Value = 6408.2
$ws.Range("J141").Value = 5681.6665
$ws.Range("K141").Value = 19224.6
$ws.Range("L141").Value = 17044.9995
$ws.Range("M141").Value = -14044.6
$ws.Range("N141").Value = -27404.9995

$ws = $wb.Worksheets.Item("ARM")
# Hunk 7 (diff @@ -9860,25 +9860,25 @@)
$ws.Range("H45").Value = 51305.953
$ws.Range("I45").Value = 83298
$ws.Range("J45").Value = 4258.8237
$ws.Range("K45").Value = 83298
$ws.Range("L45").Value = 4258.8237
$ws.Range("M45").Value = -82921
$ws.Range("N45").Value = -5012.8237

# Hunk 8 (diff @@ -10632,22 +10632,22 @@)
$ws.Range("H61").Value = 8457.532999999999
$ws.Range("I61").Value = 9797.380999999999
$ws.Range("K61").Value = 9797.380999999999
$ws.Range("M61").Value = -9585.380999999999

# Hunk 9 (diff @@ -11266,25 +11266,25 @@)
$ws.Range("H74").Value = 10575.125
$ws.Range("I74").Value = 14443.096
$ws.Range("J74").Value = 3190.818
$ws.Range("K74").Value = 14443.096
$ws.Range("L74").Value = 3190.818
$ws.Range("M74").Value = -13569.096
$ws.Range("N74").Value = -4938.818

# Hunk 10 (diff @@ -11413,25 +11413,25 @@)
$ws.Range("H77").Value = 10575.125
$ws.Range("I77").Value = 14443.096
$ws.Range("J77").Value = 3190.818
$ws.Range("K77").Value = 72215.48
$ws.Range("L77").Value = 15954.09
$ws.Range("M77").Value = -67847.48
$ws.Range("N77").Value = -24690.09

# Hunk 11 (diff @@ -12390,22 +12390,22 @@)
$ws.Range("H97").Value = 8338066
$ws.Range("I97").Value = 7142.067
$ws.Range("K97").Value = 7142.067
$ws.Range("M97").Value = -6646.067

# Hunk 12 (diff @@ -13021,22 +13021,22 @@)
$ws.Range("H110").Value = 2529.7058
$ws.Range("I110").Value = 1885
$ws.Range("K110").Value = 1885
$ws.Range("M110").Value = 160

# Hunk 13 (diff @@ -14075,22 +14075,22 @@)
$ws.Range("H132").Value = 2623.923
$ws.Range("I132").Value = 2308.2896
$ws.Range("K132").Value = 6924.8688
$ws.Range("M132").Value = -4394.8688

# Hunk 14 (diff @@ -14268,22 +14268,22 @@)
$ws.Range("H136").Value = 8457.532999999999
$ws.Range("I136").Value = 9797.380999999999
$ws.Range("K136").Value = 29392.143
$ws.Range("M136").Value = -26842.143

$ws = $wb.Worksheets.Item("BSM")
# Hunk 15 (diff @@ -18715,25 +18715,25 @@)
$ws.Range("H86").Value = 4453.778
$ws.Range("I86").Value = 7807.4614
$ws.Range("J86").Value = 1339.6428
$ws.Range("K86").Value = 7807.4614
$ws.Range("L86").Value = 1339.6428
$ws.Range("M86").Value = -6684.4614
$ws.Range("N86").Value = -3585.6428

# Hunk 16 (diff @@ -18868,25 +18868,25 @@)
$ws.Range("H89").Value = 4453.778
$ws.Range("I89").Value = 7807.4614
$ws.Range("J89").Value = 1339.6428
$ws.Range("K89").Value = 39037.307
$ws.Range("L89").Value = 6698.214
$ws.Range("M89").Value = -33421.307
$ws.Range("N89").Value = -17930.214

# Hunk 17 (diff @@ -19661,22 +19661,22 @@)
$ws.Range("H105").Value = 132401.12
$ws.Range("I105").Value = 254302.25
$ws.Range("K105").Value = 254302.25
$ws.Range("M105").Value = -252555.25

# Hunk 18 (diff @@ -21046,25 +21046,25 @@)
$ws.Range("H134").Value = 6033.6294
$ws.Range("I134").Value = 6073.385
$ws.Range("J134").Value = 5000
$ws.Range("K134").Value = 18220.155
$ws.Range("L134").Value = 15000
$ws.Range("M134").Value = -15685.155
$ws.Range("N134").Value = -20070

$ws = $wb.Worksheets.Item("CRP")
# Hunk 19 (diff @@ -24270,22 +24270,22 @@)
$ws.Range("H58").Value = 2585.8076
$ws.Range("I58").Value = 2496.3333
$ws.Range("K58").Value = 2496.3333
$ws.Range("M58").Value = -2293.3333

# Hunk 20 (diff @@ -25633,22 +25633,22 @@)
$ws.Range("H86").Value = 10235.571
$ws.Range("I86").Value = 7930.2
$ws.Range("K86").Value = 7930.2
$ws.Range("M86").Value = -6807.2

# Hunk 21 (diff @@ -25786,22 +25786,22 @@)
$ws.Range("H89").Value = 10235.571
$ws.Range("I89").Value = 7930.2
$ws.Range("K89").Value = 39651
$ws.Range("M89").Value = -34035

# Hunk 22 (diff @@ -27988,22 +27988,22 @@)
$ws.Range("H134").Value = 5197.4287
$ws.Range("I134").Value = 6756.5654
$ws.Range("K134").Value = 20269.6962
$ws.Range("M134").Value = -17734.6962

# Hunk 23 (diff @@ -28089,22 +28089,22 @@)
$ws.Range("H136").Value = 2585.8076
$ws.Range("I136").Value = 2496.3333
$ws.Range("K136").Value = 7488.999899999999
$ws.Range("M136").Value = -4938.999899999999

$ws = $wb.Worksheets.Item("CUL")
# Hunk 24 (diff @@ -30059,23 +30059,26 @@)
$ws.Range("H33").Value = 462.125
$ws.Range("J33").Value = 283.33334
$ws.Range("L33").Value = 1700.00004
$ws.Range("N33").Value = -2266.00004

# Hunk 25 (diff @@ -30108,25 +30111,25 @@)
$ws.Range("H34").Value = 2729.6667
$ws.Range("J34").Value = 3037.5
$ws.Range("L34").Value = 9112.5
$ws.Range("N34").Value = -9280.5

# Hunk 26 (diff @@ -30365,22 +30368,22 @@)
$ws.Range("H39").Value = 1024.8125
$ws.Range("I39").Value = 293.13333
$ws.Range("K39").Value = 879.39999
$ws.Range("M39").Value = -585.39999

# Hunk 27 (diff @@ -31179,25 +31182,25 @@)
$ws.Range("H55").Value = 12919.214
$ws.Range("J55").Value = 14219.36
$ws.Range("L55").Value = 42658.08
$ws.Range("N55").Value = -43012.08

$ws = $wb.Worksheets.Item("GSM")
# Hunk 28 (diff @@ -36544,7 +36547,7 @@)
$ws.Range("H20").Value = 319507

$ws = $wb.Worksheets.Item("LTW")
# Hunk 29 (diff @@ -48247,22 +48250,22 @@)
$ws.Range("H122").Value = 4950.875
$ws.Range("I122").Value = 4827.826
$ws.Range("K122").Value = 14483.478
$ws.Range("M122").Value = -12033.478

# Hunk 30 (diff @@ -48930,25 +48933,25 @@)
$ws.Range("H136").Value = 4305.0386
$ws.Range("I136").Value = 2704.2856
$ws.Range("J136").Value = 6172.5835
$ws.Range("K136").Value = 8112.8568
$ws.Range("L136").Value = 18517.7505
$ws.Range("M136").Value = -5562.8568
$ws.Range("N136").Value = -23617.7505

$ws = $wb.Worksheets.Item("WVR")
# Hunk 31 (diff @@ -50571,19 +50574,22 @@)
$ws.Range("H28").Value = 8000
$ws.Range("I28").Value = 8000
$ws.Range("K28").Value = 8000
$ws.Range("M28").Value = -7652

# Hunk 32 (diff @@ -52177,26 +52183,23 @@)
$ws.Range("H62").Value = 571439
$ws.Range("J62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("N62").ClearContents()

# Hunk 33 (diff @@ -52321,26 +52324,23 @@)
$ws.Range("H65").Value = 571439
$ws.Range("J65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("N65").ClearContents()

# Hunk 34 (diff @@ -53090,25 +53090,25 @@)
$ws.Range("H81").Value = 8567.267
$ws.Range("I81").Value = 15283.857
$ws.Range("J81").Value = 2690.25
$ws.Range("K81").Value = 30567.714
$ws.Range("L81").Value = 5380.5
$ws.Range("M81").Value = -29506.714
$ws.Range("N81").Value = -7502.5

# Hunk 35 (diff @@ -53237,25 +53237,25 @@)
$ws.Range("H84").Value = 8567.267
$ws.Range("I84").Value = 15283.857
$ws.Range("J84").Value = 2690.25
$ws.Range("K84").Value = 152838.57
$ws.Range("L84").Value = 26902.5
$ws.Range("M84").Value = -147534.57
$ws.Range("N84").Value = -37510.5

# Hunk 36 (diff @@ -54839,22 +54839,22 @@)
$ws.Range("H117").Value = 68409
$ws.Range("J117").Value = 68409
$ws.Range("L117").Value = 68409
$ws.Range("N117").Value = -77587

# Hunk 37 (diff @@ -55758,25 +55758,25 @@)
$ws.Range("H136").Value = 654738.75
$ws.Range("I136").Value = 916537.0600000001
$ws.Range("J136").Value = 18942.857
$ws.Range("K136").Value = 2749611.18
$ws.Range("L136").Value = 56828.571
$ws.Range("M136").Value = -2747061.18
$ws.Range("N136").Value = -61928.571
